# Natmi following Dr Hou advice
# Update the Cxcl13-Ccr10 LR-pair sheet: add ECs (Endothelial Cells) rows
# and recompute sending/target cluster combinations (ECs/FAPs x FAPs/sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2087556666666667
$ws.Range("H2").Value = 0.626267
$ws.Range("I2").Value = 0.01876624903294638
$ws.Range("J2").Value = 0.01876624903294638
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.036595333333333
$ws.Range("N2").Value = 3.109786
$ws.Range("O2").Value = 0.393072250513715
$ws.Range("P2").Value = 0.393072250513715
$ws.Range("Q2").Value = 0.2163951498735556
$ws.Range("R2").Value = 1.947556348862
$ws.Range("S2").Value = 0.007376491741081061
$ws.Range("T2").Value = 0.007376491741081061

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2087556666666667
$ws.Range("H3").Value = 0.626267
$ws.Range("I3").Value = 0.01876624903294638
$ws.Range("J3").Value = 0.01876624903294638
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.600567
$ws.Range("N3").Value = 4.801701
$ws.Range("O3").Value = 0.6069277494862849
$ws.Range("P3").Value = 0.6069277494862849
$ws.Range("Q3").Value = 0.3341274311296667
$ws.Range("R3").Value = 3.007146880167
$ws.Range("S3").Value = 0.01138975729186532
$ws.Range("T3").Value = 0.01138975729186532

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.91523966666667
$ws.Range("H4").Value = 32.745719
$ws.Range("I4").Value = 0.9812337509670536
$ws.Range("J4").Value = 0.9812337509670537
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.036595333333333
$ws.Range("N4").Value = 3.109786
$ws.Range("O4").Value = 0.393072250513715
$ws.Range("P4").Value = 0.393072250513715
$ws.Range("Q4").Value = 11.31468650068155
$ws.Range("R4").Value = 101.832178506134
$ws.Range("S4").Value = 0.385695758772634
$ws.Range("T4").Value = 0.385695758772634

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.91523966666667
$ws.Range("H5").Value = 32.745719
$ws.Range("I5").Value = 0.9812337509670536
$ws.Range("J5").Value = 0.9812337509670537
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.600567
$ws.Range("N5").Value = 4.801701
$ws.Range("O5").Value = 0.6069277494862849
$ws.Range("P5").Value = 0.6069277494862849
$ws.Range("Q5").Value = 17.47057240755767
$ws.Range("R5").Value = 157.235151668019
$ws.Range("S5").Value = 0.5955379921944196
$ws.Range("T5").Value = 0.5955379921944196
